# Update "想去人数" (F column) counts as captured in the generated-output commit.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 3445
$ws1.Range("F5").Value  = 8064
$ws1.Range("F8").Value  = 2048
$ws1.Range("F12").Value = 526
$ws1.Range("F15").Value = 1038
$ws1.Range("F17").Value = 143
$ws1.Range("F18").Value = 1108
$ws1.Range("F19").Value = 692
$ws1.Range("F20").Value = 509
$ws1.Range("F22").Value = 404
$ws1.Range("F23").Value = 1384
$ws1.Range("F24").Value = 4347
$ws1.Range("F25").Value = 76
$ws1.Range("F26").Value = 44564
$ws1.Range("F27").Value = 3844
$ws1.Range("F29").Value = 976
$ws1.Range("F30").Value = 677
$ws1.Range("F32").Value = 808
$ws1.Range("F35").Value = 177
$ws1.Range("F37").Value = 555
$ws1.Range("F38").Value = 449
$ws1.Range("F39").Value = 878
$ws1.Range("F40").Value = 102
$ws1.Range("F41").Value = 137
$ws1.Range("F42").Value = 1036
$ws1.Range("F43").Value = 660
$ws1.Range("F45").Value = 45
$ws1.Range("F47").Value = 2436

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 101
$ws2.Range("F13").Value = 38
$ws2.Range("F17").Value = 20
$ws2.Range("F18").Value = 140
$ws2.Range("F19").Value = 7219

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value  = 2179
$ws3.Range("F5").Value  = 1457
$ws3.Range("F8").Value  = 2302
$ws3.Range("F9").Value  = 9201
$ws3.Range("F10").Value = 1463
$ws3.Range("F11").Value = 138
$ws3.Range("F12").Value = 48

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 3445
$ws4.Range("F4").Value  = 2179
$ws4.Range("F6").Value  = 8064
$ws4.Range("F8").Value  = 1463
$ws4.Range("F9").Value  = 138
$ws4.Range("F10").Value = 48
$ws4.Range("F12").Value = 526
$ws4.Range("F14").Value = 1038
$ws4.Range("F18").Value = 143
$ws4.Range("F19").Value = 1108
$ws4.Range("F20").Value = 692
$ws4.Range("F22").Value = 4347
$ws4.Range("F23").Value = 76
$ws4.Range("F24").Value = 101
$ws4.Range("F26").Value = 38
$ws4.Range("F27").Value = 3845
$ws4.Range("F29").Value = 677
$ws4.Range("F31").Value = 808
$ws4.Range("F34").Value = 177
$ws4.Range("F35").Value = 20
$ws4.Range("F36").Value = 449
$ws4.Range("F37").Value = 878
$ws4.Range("F38").Value = 102
$ws4.Range("F39").Value = 137
$ws4.Range("F40").Value = 1036
$ws4.Range("F42").Value = 660
$ws4.Range("F47").Value = 2436
